$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segment names currently stored in column A (rows 2-20), in order
$segments = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

$lastRow = $segments.Length + 1   # row 20

# Insert a new column before column B, shifting old B:E -> C:F
$ws.Columns("B").Insert()

# The inserted column's header (B1) lost the bold/border header style, and
# the inserted column's data cells (B2:B20) incorrectly inherited column A's
# bold/border style. Fix both up to match the target layout:
#  - B1 ("segments" header) should carry the same header style as the other
#    header cells (copy format from C1, which kept it).
#  - B2:B20 (segment-name cells) should be unstyled, like the rest of the
#    data cells.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range($ws.Cells.Item(2, 2), $ws.Cells.Item($lastRow, 2)).ClearFormats()

# New header text in B1
$ws.Cells.Item(1, 2).Value = "segments"

# Fill column A with a 0-based numeric index and column B with the segment
# names (previously held in column A before the insert).
for ($i = 0; $i -lt $segments.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $segments[$i]
}
